# Auto update stock data
# - Bump the "Date_1" column (A) from 2025/11/29 -> 2025/11/30 for every
#   company's latest-data row.
# - Fill in the previously-blank Altman Z-Score / Piotroski F-Score /
#   Risk Level columns (G/H/J) for the "Ultra" company block (rows 44-49).
# - Refresh MKS's EBITDA value (B74) from 16.06 -> 16.17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to stay a plain text
# string (column A holds dates formatted as literal "yyyy/mm/dd" text, not
# real date serials, and the numeric-looking EBITDA column is text too) -
# without leaving the cell's number format/style changed afterwards.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Rows whose Date_1 cell (column A) needs to move from 2025/11/29 to
# 2025/11/30 (the most-recent-data row for each company in the table).
$dateRows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $dateRows) {
    Set-TextValue $ws.Cells.Item($r, 1) "2025/11/30"
}

# Ultra (rows 44-49): populate Altman Z-Score, Piotroski F-Score and the
# derived Risk Level, which were previously left blank.
foreach ($r in 44..49) {
    $ws.Cells.Item($r, 7).Value = 2.77
    $ws.Cells.Item($r, 8).Value = 6
    $ws.Cells.Item($r, 10).Value = "Medium risk"
}

# MKS (row 74): refreshed EBITDA figure.
Set-TextValue $ws.Cells.Item(74, 2) "16.17"
